# Add source projects/clients to UNDER exceptions, and shift employee
# names/details to reflect the current allocation period per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the StartDate/EndDate columns (C, D) store plain text dates
# (e.g. "2026-01-18"), not real Excel dates. Force text format on a
# cell before writing a date-like string so Excel keeps it as text
# instead of auto-converting it to a date serial number.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Chen Noam: add source
$ws.Range("F2").Value = "Sightec (Israel) Ltd."

# Row 4 - Meir Zipora: add source
$ws.Range("F4").Value = "BioXtreme Ltd."

# Row 6 - now Weingarten Ayala (was Noiman Yehuda)
$ws.Range("A6").Value = "Weingarten Ayala"
Set-TextValue $ws.Range("D6") "2026-02-28"
$ws.Range("E6").Value = 85
$ws.Range("F6").Value = "Arad Technologies Ltd."

# Row 7 - now Itzhaki Yair (was Weingarten Ayala)
$ws.Range("A7").Value = "Itzhaki Yair"
Set-TextValue $ws.Range("D7") "2026-04-18"
$ws.Range("E7").Value = 50

# Row 8 - now Peretz Yehonathan (was Itzhaki Yair)
$ws.Range("A8").Value = "Peretz Yehonathan"
Set-TextValue $ws.Range("D8") "2026-03-31"
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = "Red Sea"

# Row 9 - now Levin Yanir (was Peretz Yehonathan)
$ws.Range("A9").Value = "Levin Yanir"
Set-TextValue $ws.Range("D9") "2026-02-28"
$ws.Range("F9").Value = "Speedata Ltd"

# Row 10 - now Halevy Maor (was Levin Yanir)
$ws.Range("A10").Value = "Halevy Maor"
$ws.Range("F10").Value = "Rav Bariach Locks Products LTD."

# Row 11 - now Cohen Aharon (was Halevy Maor)
$ws.Range("A11").Value = "Cohen Aharon"
Set-TextValue $ws.Range("D11") "2026-03-31"
$ws.Range("F11").Value = "Aquestia Ltd."

# Row 12 - now Pruzanski Yossi (was Cohen Aharon)
$ws.Range("A12").Value = "Pruzanski Yossi"
$ws.Range("F12").Value = "Maytronics Ltd."

# Row 13 - now Morgenstern Elisheva (was Pruzanski Yossi)
$ws.Range("A13").Value = "Morgenstern Elisheva"
$ws.Range("F13").Value = "Red Sea"

# Row 14 - now Noiman Yehuda (was Morgenstern Elisheva)
$ws.Range("A14").Value = "Noiman Yehuda"
Set-TextValue $ws.Range("C14") "2026-02-01"
Set-TextValue $ws.Range("D14") "2026-04-18"
